$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $oldStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $oldStyle
}

$ws.Range("D2").Value = "71.016.47"
$ws.Range("E2").Value = "  +2.77%  "
$ws.Range("D3").Value = "3.788.93"
$ws.Range("E3").Value = "  +0.51%  "
$ws.Range("E4").Value = "  +0.00%  "
Set-TextValue "D5" "703.25"
$ws.Range("E5").Value = "  +11.17%  "
Set-TextValue "D6" "173.01"
$ws.Range("E6").Value = "  +4.00%  "
$ws.Range("D7").Value = "3.787.41"
$ws.Range("E7").Value = "  +0.55%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  +0.84%  "
Set-TextValue "D10" "0.162"
Set-TextValue "D11" "7.44"
$ws.Range("E11").Value = "  +9.53%  "
$ws.Range("E12").Value = "  +0.45%  "
$ws.Range("E13").Value = "  +6.54%  "
$ws.Range("E14").Value = "  +3.54%  "
$ws.Range("D15").Value = "4.427.11"
$ws.Range("E15").Value = "  +0.52%  "
$ws.Range("D16").Value = "3.790.23"
$ws.Range("E16").Value = "  +1.15%  "
$ws.Range("D17").Value = "71.030.85"
$ws.Range("E17").Value = "  +2.79%  "
$ws.Range("E18").Value = "  +1.37%  "
$ws.Range("E19").Value = "  +2.39%  "
$ws.Range("E20").Value = "  +0.85%  "
$ws.Range("E21").Value = "  +17.72%  "
Set-TextValue "D22" "481.36"
$ws.Range("E22").Value = "  +4.21%  "
$ws.Range("E23").Value = "  +1.21%  "
Set-TextValue "D24" "83.90"
$ws.Range("E25").Value = "  +0.25%  "
Set-TextValue "D26" "12.37"
$ws.Range("E26").Value = "  +2.10%  "
Set-TextValue "D27" "2.18"
$ws.Range("E27").Value = "  +2.27%  "
Set-TextValue "D28" "10.54"
$ws.Range("E28").Value = "  +4.63%  "
$ws.Range("D29").Value = "3.939.38"
$ws.Range("E29").Value = "  +0.52%  "
Set-TextValue "D30" "1.00"
$ws.Range("E30").Value = "  -0.08%  "
$ws.Range("E31").Value = "  +15.67%  "
$ws.Range("E32").Value = "  +0.97%  "
Set-TextValue "D33" "7.56"
$ws.Range("E33").Value = "  +7.07%  "
Set-TextValue "D34" "29.53"
$ws.Range("E34").Value = "  +3.91%  "
Set-TextValue "D35" "0.177"
$ws.Range("E35").Value = "  +0.28%  "
$ws.Range("E36").Value = "  +2.78%  "
$ws.Range("E37").Value = "  +0.15%  "
$ws.Range("D38").Value = "3.739.13"
$ws.Range("E38").Value = "  +0.42%  "
$ws.Range("E39").Value = "  +1.94%  "
Set-TextValue "D40" "3.46"
$ws.Range("E40").Value = "  +4.94%  "
$ws.Range("E41").Value = "  +3.10%  "
Set-TextValue "D42" "2.24"
$ws.Range("E42").Value = "  +13.43%  "
$ws.Range("E43").Value = "  +22.51%  "
Set-TextValue "D44" "0.966"
$ws.Range("E44").Value = "  +0.39%  "
$ws.Range("E45").Value = "  +0.00%  "
$ws.Range("E46").Value = "  +0.01%  "
Set-TextValue "D47" "46.16"
$ws.Range("E47").Value = "  +7.31%  "
Set-TextValue "D48" "160.92"
$ws.Range("E48").Value = "  +2.22%  "
$ws.Range("B49").Value = "ONDO"
$ws.Range("C49").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
Set-TextValue "D49" "1.43"
$ws.Range("E49").Value = "  -0.44%  "
$ws.Range("B50").Value = "OKB"
$ws.Range("C50").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D50" "49.07"
$ws.Range("E50").Value = "  +4.50%  "
$ws.Range("E51").Value = "  +1.81%  "
